# Adds an "alias" column to the "Experiências" and "Formações" sheets,
# right before the existing "cargo"/"curso" columns, and fills it with a
# short alias/slug for each company / institution already listed in that
# row. Also updates the active sheet/selection to reflect where the user
# ended up working (Experiências, cell C1) and leaves Formações selected
# at B7 from inserting its new column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Experiências": insert new column C ("alias") before "cargo_pt"
# ---------------------------------------------------------------------
$wsExp = $wb.Worksheets.Item("Experiências")
$wsExp.Columns.Item(3).Insert()

$wsExp.Cells.Item(1, 3).Value2 = "alias"

# Fill row by row in the same order the workbook's author typed them in
# (this also controls the order new shared strings are appended in).
$wsExp.Cells.Item(2, 3).Value2 = "nov"   # Novid
$wsExp.Cells.Item(5, 3).Value2 = "eco"   # EcoEsfera
$wsExp.Cells.Item(6, 3).Value2 = "yop"   # Yoppi
$wsExp.Cells.Item(3, 3).Value2 = "luna"  # LunaTech
$wsExp.Cells.Item(4, 3).Value2 = "aqn"   # Aquino & Silva
$wsExp.Cells.Item(7, 3).Value2 = "pipa"  # PIPA/UFLB

# ---------------------------------------------------------------------
# Sheet "Formações": insert new column B ("alias") before "curso_pt"
# ---------------------------------------------------------------------
$wsForm = $wb.Worksheets.Item("Formações")
$wsForm.Columns.Item(2).Insert()

$wsForm.Cells.Item(1, 2).Value2 = "alias"

$wsForm.Cells.Item(2, 2).Value2 = "grad"  # Universidade Federal Líbero Badaró
$wsForm.Cells.Item(3, 2).Value2 = "mic"   # MicroAdvance
$wsForm.Cells.Item(4, 2).Value2 = "fer"   # FERVESP
$wsForm.Cells.Item(5, 2).Value2 = "csap"  # CSAP
$wsForm.Cells.Item(6, 2).Value2 = "agil"  # AgileBuddies

# Match the new "alias" column's width to the rest of the sheet.
$wsForm.Columns.Item(2).ColumnWidth = 14.333333333333334

# ---------------------------------------------------------------------
# Selection / active sheet bookkeeping, matching the saved view state.
# ---------------------------------------------------------------------
$wsForm.Activate()
$wsForm.Range("B7").Select()

$wsExp.Activate()
$wsExp.Range("C1").Select()
